$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up player name spellings (shared-string text corrections)
$ws.Range("E3").Value = "Erik van Rooyen"
$ws.Range("D7").Value = "Rafael Cabrera Bello"
$ws.Range("E9").Value = "Alexander Noren"
$ws.Range("E11").Value = "Romain Langasque"

# Set explicit column widths (values chosen so the stored width - which is
# quantized to 1/6-character-unit increments - lands as close as possible to
# the target widths 19.1640625 / 23.33203125 / 27.33203125 / 26.1640625 / 31.6640625)
$ws.Columns.Item(1).ColumnWidth = 18.3333333333333
$ws.Columns.Item(2).ColumnWidth = 22.5
$ws.Columns.Item(3).ColumnWidth = 26.5
$ws.Columns.Item(4).ColumnWidth = 25.3333333333333
$ws.Columns.Item(5).ColumnWidth = 30.8333333333333

# Move selection to E12
$ws.Range("E12").Select()
